$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alumni_seminars")

# --- Step 1: clone formatting from template rows (145 for rows with G/H; 147 for rows without) ---
$ws.Range("A145:I145").Copy() | Out-Null
$ws.Range("A148:I148").PasteSpecial(-4122) | Out-Null
$ws.Range("A145:I145").Copy() | Out-Null
$ws.Range("A149:I149").PasteSpecial(-4122) | Out-Null
$ws.Range("A145:I145").Copy() | Out-Null
$ws.Range("A150:I150").PasteSpecial(-4122) | Out-Null
$ws.Range("A147:I147").Copy() | Out-Null
$ws.Range("A151:I151").PasteSpecial(-4122) | Out-Null
$ws.Range("A147:I147").Copy() | Out-Null
$ws.Range("A152:I152").PasteSpecial(-4122) | Out-Null
$ws.Range("A147:I147").Copy() | Out-Null
$ws.Range("A153:I153").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- Step 2: clear F (Slide) cells picked up from column default style; not present in target rows ---
$ws.Cells.Item(148,6).Clear() | Out-Null
$ws.Cells.Item(149,6).Clear() | Out-Null
$ws.Cells.Item(150,6).Clear() | Out-Null
$ws.Cells.Item(151,6).Clear() | Out-Null
$ws.Cells.Item(152,6).Clear() | Out-Null
$ws.Cells.Item(153,6).Clear() | Out-Null

# --- Step 3: clear G/H on rows 151-153 (not present in target; template row 147 had none, but paste still created default-styled empties) ---
$ws.Cells.Item(151,7).Clear() | Out-Null
$ws.Cells.Item(151,8).Clear() | Out-Null
$ws.Cells.Item(152,7).Clear() | Out-Null
$ws.Cells.Item(152,8).Clear() | Out-Null
$ws.Cells.Item(153,7).Clear() | Out-Null
$ws.Cells.Item(153,8).Clear() | Out-Null

# --- Step 4: apply J-column style (s=8) to rows 148 and 152 by copying format from J114 ---
$ws.Range("J114").Copy() | Out-Null
$ws.Range("J148").PasteSpecial(-4122) | Out-Null
$ws.Range("J152").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- Step 5: date values (column A) ---
$ws.Cells.Item(148,1).Value = 45499
$ws.Cells.Item(149,1).Value = 45513
$ws.Cells.Item(150,1).Value = 45527
$ws.Cells.Item(151,1).Value = 45541
$ws.Cells.Item(152,1).Value = 45549
$ws.Cells.Item(153,1).Value = 45555

# --- Step 6: new text values, written in the original authoring order so the shared-string table lines up ---
$ws.Range("B148").Value = "万维钢"
$ws.Range("C148").Value = "AI 奇点降至？"
$ws.Range("C149").Value = "比才和他的歌剧《卡门》"
$ws.Range("C150").Value = "Surprise Topic"
$ws.Range("D148").Value = "img/AI奇点.jpg"
$ws.Range("D149").Value = "img/Opera-Carmen.jpg"
$ws.Range("D150").Value = "img/suprise.jpg"
$ws.Range("D153").Value = "img/jazz.jpg"
$ws.Range("C151").Value = "National Association of Realtors Settlement and Impacts"
$ws.Range("D151").Value = "img/NAR-settlement-impact.jpg"
$ws.Range("E151").Value = "Upcoming"
$ws.Range("B152").Value = "李惠南"
$ws.Range("C152").Value = "老年痴呆和干细胞研究最新进展"
$ws.Range("C153").Value = "Jazz， Ultimate American Art Form 1:  an Emotion Conduit"
$ws.Range("J148").Value = "物理博士， 著名科普畅销书作家，得到APP《精英日课》专栏作家"
$ws.Range("J152").Value = "Group lead in Neuroscience, UCSF基因筛查中心， 分别师从老年痴呆研究和干细胞研究获得者"
$ws.Range("D152").Value = "img/stemcells.jpg"
$ws.Range("G148").Value = "https://drive.google.com/file/d/15E2JeMmO88Uq6BTzopYpKA53qP-Wt3wC/view?usp=sharing"
$ws.Range("H148").Value = "https://drive.google.com/file/d/1bVpg6sKSiF6s-sWR2xa5yLrzFH702mOv/view?usp=sharing"
$ws.Range("G149").Value = "https://drive.google.com/file/d/1tsBLYmPLeryg6eZnkbqjSFWZumnHLA6v/view?usp=sharing"
$ws.Range("H149").Value = "https://drive.google.com/file/d/17kWjulICpxc-TtfSNyCA2UO-0Vutg5rU/view?usp=sharing"
$ws.Range("G150").Value = "https://drive.google.com/file/d/1bAEN5lZrCAwWrvQwrCsCdt2ChYOJud5A/view?usp=sharing"
$ws.Range("H150").Value = "https://drive.google.com/file/d/12q6WI1EzH89J5bGpXVNpkATvK1v8Y2Se/view?usp=sharing"

# --- Step 7: remaining text values that reuse already-existing shared strings ---
$ws.Range("B149").Value = "吴鹏"
$ws.Range("B150").Value = "苏玻"
$ws.Range("B151").Value = "彭松石"
$ws.Range("B153").Value = "赵书来"
$ws.Range("E148").Value = "Technology"
$ws.Range("E149").Value = "Hobby"
$ws.Range("E150").Value = "Experience"
$ws.Range("E152").Value = "Upcoming"
$ws.Range("E153").Value = "Upcoming"

# --- Step 8: StartTime numeric values (column I) ---
$ws.Range("I151").Value = 0.8125
$ws.Range("I152").Value = 0.8125
$ws.Range("I153").Value = 0.8125

# --- Step 9: restore selection state to match the saved view (row 151 selected) ---
$ws.Rows(151).Select() | Out-Null

Write-Host "done"
